$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053203758824112
$ws.Range("D2").Value = 1.054741524300087
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.063515329478111
$ws.Range("I2").Value = 1.041614663662153
$ws.Range("J2").Value = 1.058222334472226
$ws.Range("K2").Value = 1.057483898999779
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.066233760073593
$ws.Range("N2").Value = 1.059725131566975

$ws.Range("B3").Value = 1.019999999999999
$ws.Range("C3").Value = 1.054933369543026
$ws.Range("D3").Value = 1.056077087040096
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.065089258232007
$ws.Range("I3").Value = 1.042068476662145
$ws.Range("J3").Value = 1.059598948930413
$ws.Range("K3").Value = 1.058631431466922
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.067620811408918
$ws.Range("N3").Value = 1.061103700975592

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056049328947736
$ws.Range("D4").Value = 1.05693826751714
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.066104926468444
$ws.Range("I4").Value = 1.042359091091859
$ws.Range("J4").Value = 1.060486197599019
$ws.Range("K4").Value = 1.059370409719171
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.068515069344284
$ws.Range("N4").Value = 1.061992209639076

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056517723996725
$ws.Range("D5").Value = 1.05729959618471
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.066531263162632
$ws.Range("I5").Value = 1.042480544033495
$ws.Range("J5").Value = 1.060858368139762
$ws.Range("K5").Value = 1.059680236059867
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.068890247393726
$ws.Range("N5").Value = 1.062364908704676

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056596325795191
$ws.Range("D6").Value = 1.057360223443494
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.066602809204937
$ws.Range("I6").Value = 1.042500894351889
$ws.Range("J6").Value = 1.060920808962365
$ws.Range("K6").Value = 1.059732208296423
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.068953196743312
$ws.Range("N6").Value = 1.062427438200409

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056055590608266
$ws.Range("D7").Value = 1.056943098389623
$ws.Range("E7").Value = 0.9943035907978918
$ws.Range("F7").Value = 1.066110625737813
$ws.Range("I7").Value = 1.042360716779956
$ws.Range("J7").Value = 1.060491173799873
$ws.Range("K7").Value = 1.059374552920917
$ws.Range("L7").Value = 0.9968970624459044
$ws.Range("M7").Value = 1.068520085493264
$ws.Range("N7").Value = 1.061997192906707

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053788963641975
$ws.Range("D8").Value = 1.055193515449165
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.064047827041119
$ws.Range("I8").Value = 1.041768662737047
$ws.Range("J8").Value = 1.058688302531076
$ws.Range("K8").Value = 1.057872454623359
$ws.Range("L8").Value = 0.9958175282591057
$ws.Range("M8").Value = 1.066703202062891
$ws.Range("N8").Value = 1.060191761353916

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04976957142443
$ws.Range("D9").Value = 1.052086921409484
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.060391134812691
$ws.Range("I9").Value = 1.0407019455836
$ws.Range("J9").Value = 1.055483957357731
$ws.Range("K9").Value = 1.055197899394579
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.063476145986216
$ws.Range("N9").Value = 1.056982865642881

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.047071996034689
$ws.Range("D10").Value = 1.049999295618122
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.057937895128154
$ws.Range("I10").Value = 1.039974737903881
$ws.Range("J10").Value = 1.0533284928222
$ws.Range("K10").Value = 1.053395610863894
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.06130690075262
$ws.Range("N10").Value = 1.054824340100474

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.045899434935521
$ws.Range("D11").Value = 1.049091246683636
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.056871771875845
$ws.Range("I11").Value = 1.039655971078841
$ws.Range("J11").Value = 1.052390416982893
$ws.Range("K11").Value = 1.052610487747038
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.060363186595165
$ws.Range("N11").Value = 1.053884932085877

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.045463199726864
$ws.Range("D12").Value = 1.048753328413184
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.056475171534662
$ws.Range("I12").Value = 1.039536977846111
$ws.Range("J12").Value = 1.052041245801429
$ws.Range("K12").Value = 1.052318136188661
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.060011971037076
$ws.Range("N12").Value = 1.053535265041285

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045556805359837
$ws.Range("D13").Value = 1.04882584168733
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.056560270818402
$ws.Range("I13").Value = 1.039562529050145
$ws.Range("J13").Value = 1.052116177395065
$ws.Range("K13").Value = 1.052380879462351
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.060087338887473
$ws.Range("N13").Value = 1.053610303046379

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.04586338985991
$ws.Range("D14").Value = 1.049063327153241
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.056839000990668
$ws.Range("I14").Value = 1.039646147112059
$ws.Range("J14").Value = 1.052361569324022
$ws.Range("K14").Value = 1.052586336676121
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.06033416892543
$ws.Range("N14").Value = 1.053856043460021

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.046052194187509
$ws.Range("D15").Value = 1.049209566102018
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.05701065657508
$ws.Range("I15").Value = 1.039697588726786
$ws.Range("J15").Value = 1.052512666462427
$ws.Range("K15").Value = 1.052712829608218
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.060486158764967
$ws.Range("N15").Value = 1.054007355173693

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.047149718791115
$ws.Range("D16").Value = 1.050059472447167
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.058008567681816
$ws.Range("I16").Value = 1.039995811131877
$ws.Range("J16").Value = 1.053390648560412
$ws.Range("K16").Value = 1.053447616367466
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.061369437661577
$ws.Range("N16").Value = 1.054886584106964

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.047836950857927
$ws.Range("D17").Value = 1.050591490894565
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.058633488076272
$ws.Range("I17").Value = 1.040181834911319
$ws.Range("J17").Value = 1.053940102438859
$ws.Range("K17").Value = 1.053907256012782
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.061922302497825
$ws.Range("N17").Value = 1.055436818272942

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.048237369449848
$ws.Range("D18").Value = 1.050901414241602
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.058997622753221
$ws.Range("I18").Value = 1.04028996537073
$ws.Range("J18").Value = 1.054260132749882
$ws.Range("K18").Value = 1.054174901468224
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.062244354051648
$ws.Range("N18").Value = 1.055757303063714

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.048373829089896
$ws.Range("D19").Value = 1.05100702366861
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.059121720841181
$ws.Range("I19").Value = 1.04032677179405
$ws.Range("J19").Value = 1.054369177781418
$ws.Range("K19").Value = 1.054266084930401
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.0623540936333
$ws.Range("N19").Value = 1.055866502951703

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047763262205623
$ws.Range("D20").Value = 1.050534451202168
$ws.Range("E20").Value = 0.9894336180355766
$ws.Range("F20").Value = 1.058566478500237
$ws.Range("I20").Value = 1.040161915054403
$ws.Range("J20").Value = 1.053881198599722
$ws.Range("K20").Value = 1.053857988139142
$ws.Range("L20").Value = 0.9929783193490043
$ws.Range("M20").Value = 1.061863029413551
$ws.Range("N20").Value = 1.055377830783597

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045773127670345
$ws.Range("D21").Value = 1.04899341106074
$ws.Range("E21").Value = 0.9882828385668255
$ws.Range("F21").Value = 1.056756938445121
$ws.Range("I21").Value = 1.039621539985527
$ws.Range("J21").Value = 1.052289327736959
$ws.Range("K21").Value = 1.052525854677649
$ws.Range("L21").Value = 0.9920501090198107
$ws.Range("M21").Value = 1.060261502467955
$ws.Range("N21").Value = 1.053783699281618

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044517827483568
$ws.Range("D22").Value = 1.048020855040584
$ws.Range("E22").Value = 0.9875604150241496
$ws.Range("F22").Value = 1.055615760181882
$ws.Range("I22").Value = 1.039278373719677
$ws.Range("J22").Value = 1.051284235479659
$ws.Range("K22").Value = 1.051684106543853
$ws.Range("L22").Value = 0.991467000034148
$ws.Range("M22").Value = 1.059250628293447
$ws.Range("N22").Value = 1.052777179678027

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.045183673577905
$ws.Range("D23").Value = 1.048536775248332
$ws.Range("E23").Value = 0.9879432794636459
$ws.Range("F23").Value = 1.056221052527919
$ws.Range("I23").Value = 1.039460617981236
$ws.Range("J23").Value = 1.051817458983201
$ws.Range("K23").Value = 1.05213073393391
$ws.Range("L23").Value = 0.9917760702887607
$ws.Range("M23").Value = 1.059786889300003
$ws.Range("N23").Value = 1.053311160420107

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047796560271026
$ws.Range("D24").Value = 1.05056022620126
$ws.Range("E24").Value = 0.9894529299347241
$ws.Range("F24").Value = 1.058596758388887
$ws.Range("I24").Value = 1.040170917136078
$ws.Range("J24").Value = 1.053907816120544
$ws.Range("K24").Value = 1.053880251574979
$ws.Range("L24").Value = 0.9929938892766438
$ws.Range("M24").Value = 1.061889813682476
$ws.Range("N24").Value = 1.055404486104352

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.050811778352038
$ws.Range("D25").Value = 1.052892918546142
$ws.Range("E25").Value = 0.9912096547607046
$ws.Range("F25").Value = 1.061339140234389
$ws.Range("I25").Value = 1.040980526378144
$ws.Range("J25").Value = 1.056315690240146
$ws.Range("K25").Value = 1.055892681108022
$ws.Range("L25").Value = 0.9944092447426411
$ws.Range("M25").Value = 1.064313512592895
$ws.Range("N25").Value = 1.05781577968139

